# Avanzamento.xlsx update (lbianco via Streamlit):
# refresh the "Produzione" (D) figures for the technicians in the second
# "45930" batch (rows 66-127). Column F ("Avanzamento") is a shared formula
# (=D-(D*E)/100) so it recalculates on its own once D changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D66").Value = 30.4796224427736
$ws.Range("D67").Value = 43.670197075095402
$ws.Range("D68").Value = 26.9223577910596
$ws.Range("D69").Value = 52.653562131096599
$ws.Range("D71").Value = 38.288378542549701
$ws.Range("D73").Value = 32.274288017826002
$ws.Range("D74").Value = 61.155298750225697
$ws.Range("D75").Value = 62.703213748923702
$ws.Range("D76").Value = 53.7327838143786
$ws.Range("D82").Value = 27.1032531964238
$ws.Range("D83").Value = 31.424534920634901
$ws.Range("D85").Value = 31.4199073391813
$ws.Range("D87").Value = 40.950688020746199
$ws.Range("D88").Value = 34.866680198944998
$ws.Range("D90").Value = 93.189302289873396
$ws.Range("D92").Value = 25.3334939885636
$ws.Range("D94").Value = 44.4334738728435
$ws.Range("D96").Value = 43.805567710536202
$ws.Range("D100").Value = 41.742409790640401
$ws.Range("D101").Value = 94.654398077086
$ws.Range("D105").Value = 37.3849826388889
$ws.Range("D106").Value = 74.312100629670397
$ws.Range("D107").Value = 27.5479324259333
$ws.Range("D115").Value = 33.525773090277802
$ws.Range("D117").Value = 33.412511168981503
$ws.Range("D121").Value = 69.058091459141295
$ws.Range("D122").Value = 42.576660968486202
$ws.Range("D123").Value = 43.319398809523797
$ws.Range("D124").Value = 30.457376481564602
$ws.Range("D125").Value = 52.236506004014601
$ws.Range("D126").Value = 70.974174818238893
$ws.Range("D127").Value = 35.704538261217998

# Leave the workbook with the same cell selected/visible as in the saved file.
$ws.Range("N58").Select()
